$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna1"
$ws.Cells.Item(2, 3).Value = "Epha2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 30.194115
$ws.Cells.Item(2, 8).Value = 90.582345
$ws.Cells.Item(2, 9).Value = 0.9018420607989291
$ws.Cells.Item(2, 10).Value = 0.901842060798929
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 20.66830833333333
$ws.Cells.Item(2, 14).Value = 62.004925
$ws.Cells.Item(2, 15).Value = 0.6755285375771634
$ws.Cells.Item(2, 16).Value = 0.6755285375771634
$ws.Cells.Item(2, 17).Value = 624.061278672125
$ws.Cells.Item(2, 18).Value = 5616.551508049125
$ws.Cells.Item(2, 19).Value = 0.6092200484570759
$ws.Cells.Item(2, 20).Value = 0.6092200484570758

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna1"
$ws.Cells.Item(3, 3).Value = "Epha2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 30.194115
$ws.Cells.Item(3, 8).Value = 90.582345
$ws.Cells.Item(3, 9).Value = 0.9018420607989291
$ws.Cells.Item(3, 10).Value = 0.901842060798929
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.156330666666667
$ws.Cells.Item(3, 14).Value = 6.468992
$ws.Cells.Item(3, 15).Value = 0.07047809033489469
$ws.Cells.Item(3, 16).Value = 0.07047809033489467
$ws.Cells.Item(3, 17).Value = 65.10849612736
$ws.Cells.Item(3, 18).Value = 585.97646514624
$ws.Cells.Item(3, 19).Value = 0.06356010622879452
$ws.Cells.Item(3, 20).Value = 0.0635601062287945

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna1"
$ws.Cells.Item(4, 3).Value = "Epha2"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 30.194115
$ws.Cells.Item(4, 8).Value = 90.582345
$ws.Cells.Item(4, 9).Value = 0.9018420607989291
$ws.Cells.Item(4, 10).Value = 0.901842060798929
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.683564
$ws.Cells.Item(4, 14).Value = 5.050692
$ws.Cells.Item(4, 15).Value = 0.05502605769642779
$ws.Cells.Item(4, 16).Value = 0.05502605769642779
$ws.Cells.Item(4, 17).Value = 50.83372502586
$ws.Cells.Item(4, 18).Value = 457.50352523274
$ws.Cells.Item(4, 19).Value = 0.04962481327058722
$ws.Cells.Item(4, 20).Value = 0.04962481327058721

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efna1"
$ws.Cells.Item(5, 3).Value = "Epha2"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 30.194115
$ws.Cells.Item(5, 8).Value = 90.582345
$ws.Cells.Item(5, 9).Value = 0.9018420607989291
$ws.Cells.Item(5, 10).Value = 0.901842060798929
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.278649666666666
$ws.Cells.Item(5, 14).Value = 15.835949
$ws.Cells.Item(5, 15).Value = 0.1725288026574751
$ws.Cells.Item(5, 16).Value = 0.1725288026574751
$ws.Cells.Item(5, 17).Value = 159.384155080045
$ws.Cells.Item(5, 18).Value = 1434.457395720405
$ws.Cells.Item(5, 19).Value = 0.1555937309357891
$ws.Cells.Item(5, 20).Value = 0.155593730935789

$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Efna1"
$ws.Cells.Item(6, 3).Value = "Epha2"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 30.194115
$ws.Cells.Item(6, 8).Value = 90.582345
$ws.Cells.Item(6, 9).Value = 0.9018420607989291
$ws.Cells.Item(6, 10).Value = 0.901842060798929
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8089063333333334
$ws.Cells.Item(6, 14).Value = 2.426719
$ws.Cells.Item(6, 15).Value = 0.02643851173403914
$ws.Cells.Item(6, 16).Value = 0.02643851173403913
$ws.Cells.Item(6, 17).Value = 24.424210852895
$ws.Cells.Item(6, 18).Value = 219.817897676055
$ws.Cells.Item(6, 19).Value = 0.02384336190668252
$ws.Cells.Item(6, 20).Value = 0.02384336190668252

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna1"
$ws.Cells.Item(7, 3).Value = "Epha2"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.633202
$ws.Cells.Item(7, 8).Value = 7.899606
$ws.Cells.Item(7, 9).Value = 0.07864884657754871
$ws.Cells.Item(7, 10).Value = 0.07864884657754868
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 20.66830833333333
$ws.Cells.Item(7, 14).Value = 62.004925
$ws.Cells.Item(7, 15).Value = 0.6755285375771634
$ws.Cells.Item(7, 16).Value = 0.6755285375771634
$ws.Cells.Item(7, 17).Value = 54.42383083995
$ws.Cells.Item(7, 18).Value = 489.81447755955
$ws.Cells.Item(7, 19).Value = 0.05312954031066217
$ws.Cells.Item(7, 20).Value = 0.05312954031066215

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efna1"
$ws.Cells.Item(8, 3).Value = "Epha2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.633202
$ws.Cells.Item(8, 8).Value = 7.899606
$ws.Cells.Item(8, 9).Value = 0.07864884657754871
$ws.Cells.Item(8, 10).Value = 0.07864884657754868
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.156330666666667
$ws.Cells.Item(8, 14).Value = 6.468992
$ws.Cells.Item(8, 15).Value = 0.07047809033489469
$ws.Cells.Item(8, 16).Value = 0.07047809033489467
$ws.Cells.Item(8, 17).Value = 5.678054224128001
$ws.Cells.Item(8, 18).Value = 51.10248801715201
$ws.Cells.Item(8, 19).Value = 0.005543020513827751
$ws.Cells.Item(8, 20).Value = 0.005543020513827748

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efna1"
$ws.Cells.Item(9, 3).Value = "Epha2"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.633202
$ws.Cells.Item(9, 8).Value = 7.899606
$ws.Cells.Item(9, 9).Value = 0.07864884657754871
$ws.Cells.Item(9, 10).Value = 0.07864884657754868
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.683564
$ws.Cells.Item(9, 14).Value = 5.050692
$ws.Cells.Item(9, 15).Value = 0.05502605769642779
$ws.Cells.Item(9, 16).Value = 0.05502605769642779
$ws.Cells.Item(9, 17).Value = 4.433164091928
$ws.Cells.Item(9, 18).Value = 39.898476827352
$ws.Cells.Item(9, 19).Value = 0.004327735969533693
$ws.Cells.Item(9, 20).Value = 0.004327735969533691

$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Efna1"
$ws.Cells.Item(10, 3).Value = "Epha2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.633202
$ws.Cells.Item(10, 8).Value = 7.899606
$ws.Cells.Item(10, 9).Value = 0.07864884657754871
$ws.Cells.Item(10, 10).Value = 0.07864884657754868
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.278649666666666
$ws.Cells.Item(10, 14).Value = 15.835949
$ws.Cells.Item(10, 15).Value = 0.1725288026574751
$ws.Cells.Item(10, 16).Value = 0.1725288026574751
$ws.Cells.Item(10, 17).Value = 13.899750859566
$ws.Cells.Item(10, 18).Value = 125.097757736094
$ws.Cells.Item(10, 19).Value = 0.01356919133041594
$ws.Cells.Item(10, 20).Value = 0.01356919133041593

$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Efna1"
$ws.Cells.Item(11, 3).Value = "Epha2"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 2.633202
$ws.Cells.Item(11, 8).Value = 7.899606
$ws.Cells.Item(11, 9).Value = 0.07864884657754871
$ws.Cells.Item(11, 10).Value = 0.07864884657754868
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.8089063333333334
$ws.Cells.Item(11, 14).Value = 2.426719
$ws.Cells.Item(11, 15).Value = 0.02643851173403914
$ws.Cells.Item(11, 16).Value = 0.02643851173403913
$ws.Cells.Item(11, 17).Value = 2.130013774746
$ws.Cells.Item(11, 18).Value = 19.170123972714
$ws.Cells.Item(11, 19).Value = 0.002079358453109166
$ws.Cells.Item(11, 20).Value = 0.002079358453109164

$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Efna1"
$ws.Cells.Item(12, 3).Value = "Epha2"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.5479553333333333
$ws.Cells.Item(12, 8).Value = 1.643866
$ws.Cells.Item(12, 9).Value = 0.01636640673320273
$ws.Cells.Item(12, 10).Value = 0.01636640673320272
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 20.66830833333333
$ws.Cells.Item(12, 14).Value = 62.004925
$ws.Cells.Item(12, 15).Value = 0.6755285375771634
$ws.Cells.Item(12, 16).Value = 0.6755285375771634
$ws.Cells.Item(12, 17).Value = 11.32530978222778
$ws.Cells.Item(12, 18).Value = 101.92778804005
$ws.Cells.Item(12, 19).Value = 0.01105597480587348
$ws.Cells.Item(12, 20).Value = 0.01105597480587348

$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Efna1"
$ws.Cells.Item(13, 3).Value = "Epha2"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.5479553333333333
$ws.Cells.Item(13, 8).Value = 1.643866
$ws.Cells.Item(13, 9).Value = 0.01636640673320273
$ws.Cells.Item(13, 10).Value = 0.01636640673320272
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.156330666666667
$ws.Cells.Item(13, 14).Value = 6.468992
$ws.Cells.Item(13, 15).Value = 0.07047809033489469
$ws.Cells.Item(13, 16).Value = 0.07047809033489467
$ws.Cells.Item(13, 17).Value = 1.181572889230222
$ws.Cells.Item(13, 18).Value = 10.634156003072
$ws.Cells.Item(13, 19).Value = 0.00115347309220029
$ws.Cells.Item(13, 20).Value = 0.00115347309220029

$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Efna1"
$ws.Cells.Item(14, 3).Value = "Epha2"
$ws.Cells.Item(14, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.5479553333333333
$ws.Cells.Item(14, 8).Value = 1.643866
$ws.Cells.Item(14, 9).Value = 0.01636640673320273
$ws.Cells.Item(14, 10).Value = 0.01636640673320272
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.683564
$ws.Cells.Item(14, 14).Value = 5.050692
$ws.Cells.Item(14, 15).Value = 0.05502605769642779
$ws.Cells.Item(14, 16).Value = 0.05502605769642779
$ws.Cells.Item(14, 17).Value = 0.9225178728079999
$ws.Cells.Item(14, 18).Value = 8.302660855272
$ws.Cells.Item(14, 19).Value = 0.0009005788411844177
$ws.Cells.Item(14, 20).Value = 0.0009005788411844174

$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Efna1"
$ws.Cells.Item(15, 3).Value = "Epha2"
$ws.Cells.Item(15, 4).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.5479553333333333
$ws.Cells.Item(15, 8).Value = 1.643866
$ws.Cells.Item(15, 9).Value = 0.01636640673320273
$ws.Cells.Item(15, 10).Value = 0.01636640673320272
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 5.278649666666666
$ws.Cells.Item(15, 14).Value = 15.835949
$ws.Cells.Item(15, 15).Value = 0.1725288026574751
$ws.Cells.Item(15, 16).Value = 0.1725288026574751
$ws.Cells.Item(15, 17).Value = 2.892464237648222
$ws.Cells.Item(15, 18).Value = 26.032178138834
$ws.Cells.Item(15, 19).Value = 0.002823676557484705
$ws.Cells.Item(15, 20).Value = 0.002823676557484704

$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Efna1"
$ws.Cells.Item(16, 3).Value = "Epha2"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.5479553333333333
$ws.Cells.Item(16, 8).Value = 1.643866
$ws.Cells.Item(16, 9).Value = 0.01636640673320273
$ws.Cells.Item(16, 10).Value = 0.01636640673320272
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.8089063333333334
$ws.Cells.Item(16, 14).Value = 2.426719
$ws.Cells.Item(16, 15).Value = 0.02643851173403914
$ws.Cells.Item(16, 16).Value = 0.02643851173403913
$ws.Cells.Item(16, 17).Value = 0.4432445395171112
$ws.Cells.Item(16, 18).Value = 3.989200855654
$ws.Cells.Item(16, 19).Value = 0.0004327034364598374
$ws.Cells.Item(16, 20).Value = 0.0004327034364598373

$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Efna1"
$ws.Cells.Item(17, 3).Value = "Epha2"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.1052186666666667
$ws.Cells.Item(17, 8).Value = 0.315656
$ws.Cells.Item(17, 9).Value = 0.00314268589031943
$ws.Cells.Item(17, 10).Value = 0.003142685890319429
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 20.66830833333333
$ws.Cells.Item(17, 14).Value = 62.004925
$ws.Cells.Item(17, 15).Value = 0.6755285375771634
$ws.Cells.Item(17, 16).Value = 0.6755285375771634
$ws.Cells.Item(17, 17).Value = 2.174691845088889
$ws.Cells.Item(17, 18).Value = 19.5722266058
$ws.Cells.Item(17, 19).Value = 0.00212297400355187
$ws.Cells.Item(17, 20).Value = 0.00212297400355187

$ws.Cells.Item(18, 1).Value = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value = "Efna1"
$ws.Cells.Item(18, 3).Value = "Epha2"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.1052186666666667
$ws.Cells.Item(18, 8).Value = 0.315656
$ws.Cells.Item(18, 9).Value = 0.00314268589031943
$ws.Cells.Item(18, 10).Value = 0.003142685890319429
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 2.156330666666667
$ws.Cells.Item(18, 14).Value = 6.468992
$ws.Cells.Item(18, 15).Value = 0.07047809033489469
$ws.Cells.Item(18, 16).Value = 0.07047809033489467
$ws.Cells.Item(18, 17).Value = 0.2268862376391111
$ws.Cells.Item(18, 18).Value = 2.041976138752
$ws.Cells.Item(18, 19).Value = 0.0002214905000721318
$ws.Cells.Item(18, 20).Value = 0.0002214905000721316

$ws.Cells.Item(19, 1).Value = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value = "Efna1"
$ws.Cells.Item(19, 3).Value = "Epha2"
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.1052186666666667
$ws.Cells.Item(19, 8).Value = 0.315656
$ws.Cells.Item(19, 9).Value = 0.00314268589031943
$ws.Cells.Item(19, 10).Value = 0.003142685890319429
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 1.683564
$ws.Cells.Item(19, 14).Value = 5.050692
$ws.Cells.Item(19, 15).Value = 0.05502605769642779
$ws.Cells.Item(19, 16).Value = 0.05502605769642779
$ws.Cells.Item(19, 17).Value = 0.177142359328
$ws.Cells.Item(19, 18).Value = 1.594281233952
$ws.Cells.Item(19, 19).Value = 0.0001729296151224665
$ws.Cells.Item(19, 20).Value = 0.0001729296151224665

$ws.Cells.Item(20, 1).Value = "Resolving-Mac"
$ws.Cells.Item(20, 2).Value = "Efna1"
$ws.Cells.Item(20, 3).Value = "Epha2"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 0.3333333333333333
$ws.Cells.Item(20, 7).Value = 0.1052186666666667
$ws.Cells.Item(20, 8).Value = 0.315656
$ws.Cells.Item(20, 9).Value = 0.00314268589031943
$ws.Cells.Item(20, 10).Value = 0.003142685890319429
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 5.278649666666666
$ws.Cells.Item(20, 14).Value = 15.835949
$ws.Cells.Item(20, 15).Value = 0.1725288026574751
$ws.Cells.Item(20, 16).Value = 0.1725288026574751
$ws.Cells.Item(20, 17).Value = 0.5554124797271111
$ws.Cells.Item(20, 18).Value = 4.998712317543999
$ws.Cells.Item(20, 19).Value = 0.0005422038337853524
$ws.Cells.Item(20, 20).Value = 0.0005422038337853522

$ws.Cells.Item(21, 1).Value = "Resolving-Mac"
$ws.Cells.Item(21, 2).Value = "Efna1"
$ws.Cells.Item(21, 3).Value = "Epha2"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = 0.3333333333333333
$ws.Cells.Item(21, 7).Value = 0.1052186666666667
$ws.Cells.Item(21, 8).Value = 0.315656
$ws.Cells.Item(21, 9).Value = 0.00314268589031943
$ws.Cells.Item(21, 10).Value = 0.003142685890319429
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 0.8089063333333334
$ws.Cells.Item(21, 14).Value = 2.426719
$ws.Cells.Item(21, 15).Value = 0.02643851173403914
$ws.Cells.Item(21, 16).Value = 0.02643851173403913
$ws.Cells.Item(21, 17).Value = 0.08511204585155556
$ws.Cells.Item(21, 18).Value = 0.7660084126640001
$ws.Cells.Item(21, 19).Value = 0.00008308793778760949
$ws.Cells.Item(21, 20).Value = 0.00008308793778760945

